# DC-Colos.xlsx update — add missing "CTU" (Chengdu, China) colo row.
#
# The generated colo list gained a new row for Chengdu (CTU) that sorts
# alphabetically right before Ashburn (IAD) in the North-America-heavy
# block starting at row 272. Inserting that single row pushes every
# following row (IAD ... YHZ) down by one, which is exactly what the
# target data shows (e.g. the old last row YHZ becomes row 333, and the
# used range grows from A1:H332 to A1:H333). No other existing row's
# data actually changes — only its row number shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above the current row 272 (IAD); this shifts
# IAD..YHZ down to rows 273..333 automatically.
$ws.Rows("272:272").Insert()

# The freshly inserted row has no formatting of its own yet — copy the
# style (bold colo code, border, centered alignment) from the cell that
# used to be directly below it (now A273, "IAD") so the new "CTU" cell
# matches every other colo-code cell in column A.
$ws.Range("A273").Copy()
$ws.Range("A272").PasteSpecial(-4122)

# Populate the new Chengdu, China row. Latitude/longitude are left blank
# (matching the source data, which has no coordinates for this colo yet).
$ws.Range("A272").Value = "CTU"
$ws.Range("B272").Value = "Chengdu, China"
$ws.Range("C272").Value = "Asia"
$ws.Range("D272").Value = "Chengdu"
$ws.Range("E272").Value = "China"
$ws.Range("F272").Value = "CN"

Write-Host "Inserted CTU row at 272; sheet now spans to row 333."
